# Sync attendance_reports: normalize the "Recorded By" column (G) so that
# the actual human recorder (dnasr281@gmail.com) is listed first, ahead of
# the automated "System" entry; when no such recorder is present, the last
# listed account is promoted ahead of the remaining middle entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) {
        continue
    }

    if ($parts -contains "dnasr281@gmail.com") {
        $rest = $parts | Where-Object { $_ -ne "dnasr281@gmail.com" }
        $newParts = @("dnasr281@gmail.com") + $rest
    } elseif ($parts.Count -gt 2) {
        $middle = $parts[1..($parts.Count - 2)]
        $newParts = @($parts[0], $parts[$parts.Count - 1]) + $middle
    } else {
        $newParts = $parts
    }

    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
